# Apply the "Add files via upload" edit:
#  - C2's emailed contact text/hyperlink label changes to a new address
#    (a brand-new shared string is created for it)
#  - E2's date moves from 2022-12-22 to 2023-01-05 (serial 44917 -> 44931)
#  - Columns C/D/E are resized (C widens, D narrows, E becomes a used column)
#  - The worksheet's saved selection becomes the single cell E3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes -----------------------------------------------
$ws.Range("C2").Value = "poojasawant1030@dbatu.ac.in"
$ws.Range("E2").Value = 44931

# --- Column width changes ----------------------------------------------
# Target XML `width` attributes are 29 (col C), 5 (col D) and 8.5546875
# (col E). Excel's ColumnWidth is quantized to whole pixels internally, so
# we use values that land in the middle of the pixel bucket that rounds to
# each target width.
$ws.Columns.Item(3).ColumnWidth = 28.17   # -> width 29
$ws.Columns.Item(4).ColumnWidth = 4.17    # -> width 5
$ws.Columns.Item(5).ColumnWidth = 7.7     # -> width ~8.55 (closest pixel bucket)

# --- Selection -----------------------------------------------------------
$ws.Range("E3").Select() | Out-Null
